$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "28-07-2022"
$ws.Range("A3").ClearFormats()
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "01-08-2022"
$ws.Range("A4").ClearFormats()
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0

# Row 5
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "04-08-2022"
$ws.Range("A5").ClearFormats()
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

# Row 6
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "08-08-2022"
$ws.Range("A6").ClearFormats()
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = 0

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "11-08-2022"
$ws.Range("A7").ClearFormats()

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "15-08-2022"
$ws.Range("A8").ClearFormats()

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "18-08-2022"
$ws.Range("A9").ClearFormats()

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "22-08-2022"
$ws.Range("A10").ClearFormats()

# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("A11").ClearFormats()

# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "29-08-2022"
$ws.Range("A12").ClearFormats()
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("H12").Value = 0

# Row 13
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "01-09-2022"
$ws.Range("A13").ClearFormats()
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0

# Row 14
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "05-09-2022"
$ws.Range("A14").ClearFormats()
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("H14").Value = 0

# Row 15
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "08-09-2022"
$ws.Range("A15").ClearFormats()

# Row 16
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "12-09-2022"
$ws.Range("A16").ClearFormats()

# Row 17
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "15-09-2022"
$ws.Range("A17").ClearFormats()

# Row 18
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "19-09-2022"
$ws.Range("A18").ClearFormats()

# Row 19
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "22-09-2022"
$ws.Range("A19").ClearFormats()

# Row 20
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("A20").ClearFormats()

# Row 21
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "29-09-2022"
$ws.Range("A21").ClearFormats()

